$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price cells are plain numeric-looking strings (e.g. "578.54") that Excel would
# otherwise auto-convert to a Number on assignment. The source data stores every Price
# cell as text (note values like "67.819.52" cannot be numbers), so force Text format
# on just those specific cells before writing the new value, to preserve their type.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '67.819.52'
$ws.Range('E2').Value = '  +1.02%  '
$ws.Range('D3').Value = '3.246.59'
$ws.Range('E3').Value = '  -0.08%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '578.54'
$ws.Range('E5').Value = '  +0.13%  '
$ws.Range('D6').Value = '182.64'
$ws.Range('E6').Value = '  +3.43%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').Value = '0.596'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('E9').Value = '  +4.10%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('D11').Value = '0.415'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '3.807.48'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('E13').Value = '  +0.44%  '
$ws.Range('D14').Value = '28.68'
$ws.Range('E14').Value = '  +2.94%  '
$ws.Range('D15').Value = '67.842.04'
$ws.Range('E15').Value = '  +1.07%  '
$ws.Range('D16').Value = '0.0000172'
$ws.Range('E16').Value = '  +2.17%  '
$ws.Range('D17').Value = '3.248.78'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').Value = '5.83'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').Value = '13.55'
$ws.Range('E19').Value = '  +1.91%  '
$ws.Range('D20').Value = '379.30'
$ws.Range('E20').Value = '  +3.14%  '
$ws.Range('D21').Value = '7.63'
$ws.Range('E21').Value = '  +2.22%  '
$ws.Range('D22').Value = '0.998'
$ws.Range('E22').Value = '  -0.18%  '
$ws.Range('D23').Value = '71.36'
$ws.Range('E23').Value = '  +1.72%  '
$ws.Range('D24').Value = '0.512'
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').Value = '0.0000119'
$ws.Range('E25').Value = '  +0.18%  '
$ws.Range('D26').Value = '9.89'
$ws.Range('E26').Value = '  +1.28%  '
$ws.Range('E27').Value = '  +2.29%  '
$ws.Range('E28').Value = '  -0.15%  '
$ws.Range('D29').Value = '1.97'
$ws.Range('E29').Value = '  +0.37%  '
$ws.Range('D30').Value = '5.67'
$ws.Range('E30').Value = '  +0.77%  '
$ws.Range('D31').Value = '22.83'
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('D33').Value = '7.04'
$ws.Range('E33').Value = '  +3.78%  '
$ws.Range('D34').Value = '1.26'
$ws.Range('E34').Value = '  +2.34%  '
$ws.Range('D35').Value = '1.57'
$ws.Range('E35').Value = '  +4.14%  '
$ws.Range('D36').Value = '162.26'
$ws.Range('E36').Value = '  -5.10%  '
$ws.Range('B37').Value = 'Mantle'
$ws.Range('C37').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D37').Value = '0.838'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.84'
$ws.Range('E38').Value = '  -0.87%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '6.70'
$ws.Range('E39').Value = '  +5.85%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').Value = '26.42'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.58'
$ws.Range('E41').Value = '  +6.66%  '
$ws.Range('D42').Value = '2.59'
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '25.47'
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '41.14'
$ws.Range('E44').Value = '  +1.85%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').Value = '346.13'
$ws.Range('E45').Value = '  +3.78%  '
$ws.Range('D46').Value = '0.0685'
$ws.Range('E46').Value = '  +2.04%  '
$ws.Range('D47').Value = '2.629.81'
$ws.Range('E47').Value = '  -2.86%  '
$ws.Range('D48').Value = '0.0283'
$ws.Range('E48').Value = '  +1.22%  '
$ws.Range('D49').Value = '0.102'
$ws.Range('E49').Value = '  -0.94%  '
$ws.Range('D50').Value = '0.991'
$ws.Range('E50').Value = '  +1.18%  '
$ws.Range('D51').Value = '6.17'
$ws.Range('E51').Value = '  +2.38%  '
